$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $replace, 2)
}

Replace-Text "2024-08-29 Thursday" "2024-08-30 Friday"

Replace-Text "170÷9=" "488÷4="
Replace-Text "575÷2=" "393÷3="
Replace-Text "399÷6=" "459÷4="
Replace-Text "413÷9=" "847÷8="
Replace-Text "258÷5=" "394÷4="
Replace-Text "154÷5=" "288÷8="
Replace-Text "425÷6=" "638÷6="
Replace-Text "505÷8=" "207÷9="
Replace-Text "699÷6=" "958÷2="
Replace-Text "137÷4=" "720÷8="
Replace-Text "499÷3=" "846÷8="
Replace-Text "685÷3=" "204÷5="
Replace-Text "844÷3=" "731÷5="
Replace-Text "581÷4=" "700÷8="
Replace-Text "978÷4=" "152÷3="
Replace-Text "785÷2=" "758÷5="
Replace-Text "603÷3=" "526÷2="
Replace-Text "270÷4=" "686÷6="
Replace-Text "492÷9=" "301÷5="
Replace-Text "508÷6=" "316÷9="
Replace-Text "855÷7=" "439÷8="
Replace-Text "946÷6=" "923÷9="
Replace-Text "420÷5=" "111÷9="
Replace-Text "120÷2=" "404÷4="
Replace-Text "146÷5=" "903÷6="

Write-Output "Done applying replacements"
